# Fixed issue with ramp constraint. Added initial states for processes and for storages(states).

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1) nodes sheet: add "initial_state" column (J) with value 0 for every node
# -------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("nodes")
$wsNodes.Range("J1").Value = "initial_state"
$wsNodes.Range("J2").Value = 0
$wsNodes.Range("J3").Value = 0
$wsNodes.Range("J4").Value = 0
$wsNodes.Range("J5").Value = 0
$wsNodes.Range("J6").Value = 0
$wsNodes.Range("J7").Value = 0
$wsNodes.Range("J2:J7").HorizontalAlignment = -4108
$wsNodes.Range("J4").Select() | Out-Null

# -------------------------------------------------------------------
# 2) processes sheet: add "initial_state" column (M) -- ngchp starts at 1
# -------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("processes")
$wsProc.Range("M1").Value = "initial_state"
$wsProc.Range("M2").Value = 1
$wsProc.Range("M3").Value = 0
$wsProc.Range("M4").Value = 0
$wsProc.Range("M5").Value = 0
$wsProc.Range("M6").Value = 0
$wsProc.Range("M7").Value = 0
$wsProc.Range("M8").Value = 0
$wsProc.Range("M2:M8").HorizontalAlignment = -4108
$wsProc.Range("M2").Select() | Out-Null

# -------------------------------------------------------------------
# 3) cf sheet: just move the selection (no value changes)
# -------------------------------------------------------------------
$wsCf = $wb.Worksheets.Item("cf")
$wsCf.Range("D2").Select() | Out-Null

# -------------------------------------------------------------------
# 4) inflow sheet: fix ramp constraint input data -- first timestep now -5
# -------------------------------------------------------------------
$wsInflow = $wb.Worksheets.Item("inflow")
$wsInflow.Range("B2").Value = -5
$wsInflow.Range("C2").Value = -5
$wsInflow.Range("D2").Value = -5
$wsInflow.Range("B3").Select() | Out-Null

# -------------------------------------------------------------------
# 5) process_topology sheet: ramp_up / ramp_down values raised to 0.5
#    for the first three processes, this becomes the active sheet/tab
# -------------------------------------------------------------------
$wsTopo = $wb.Worksheets.Item("process_topology")
$wsTopo.Range("G2").Value = 0.5
$wsTopo.Range("H2").Value = 0.5
$wsTopo.Range("G3").Value = 0.5
$wsTopo.Range("H3").Value = 0.5
$wsTopo.Range("G4").Value = 0.5
$wsTopo.Range("H4").Value = 0.5
$wsTopo.Range("G5").Value = 0.5
$wsTopo.Range("H5").Value = 0.5
$wsTopo.Range("G6").Value = 0.5
$wsTopo.Range("H6").Value = 0.5
$wsTopo.Range("J6").Select() | Out-Null
